# Regenerate merged AHB files
#
# The source workbook marks rows that are the *first* row of a new
# "Vorgang" group with a grey-filled style across the whole row (and a
# bold, non-centered font in column B), and it also clears the
# "AENDERUNG"/"ÄNDERUNG" marker text out of column L for every data row
# (column L keeps its grey fill + centered alignment, it just loses the
# bold marker font and the text value).
#
# Rows 3-9 already carry the correct target styling (pre-existing
# style indices 5/6/7 in xl/styles.xml), so we simply copy the already
# -correct formatting from those template rows/cells onto the rows
# that still need it, instead of re-deriving fonts/fills by hand
# (which would otherwise mint brand-new, duplicate style entries).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that are the first row of a "Vorgang" group: the whole row
# (A:V) needs the grey "group header" styling that row 9 already has.
$groupHeaderRows = @(13, 17, 23, 27, 34, 40, 63, 67)

# Rows whose column L only needs the "marker cleared" styling that
# L3 (through L9) already has.
$markerOnlyRows = @(14, 15, 16, 18, 19, 20, 21, 22, 24, 25, 26, 28, 29, 30, 31, 32, 33, 35, 36, 38, 39, 41, 42, 43, 104, 115)

foreach ($r in $groupHeaderRows) {
    $ws.Range("A9:V9").Copy()
    $destination = $ws.Range("A" + $r + ":V" + $r)
    $destination.PasteSpecial(-4122)
    $ws.Cells.Item($r, 12).Value = ""
}

foreach ($r in $markerOnlyRows) {
    $ws.Range("L3").Copy()
    $ws.Cells.Item($r, 12).PasteSpecial(-4122)
    $ws.Cells.Item($r, 12).Value = ""
}

$excel.CutCopyMode = $false
